# Apply the "re_rank" recomputation to the workers_rank_re_female sheet.
# Rows 3/4 swap identity (prolificid/name), and rows 10/11/12 rotate identity,
# while the realeffort (F) scores are updated with newly computed values and
# the race (G) values follow the corresponding identity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 / Row 4 swap (Jewel <-> Colleen) ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("D3").Value = "Colleen"
$ws.Range("G3").Value = "White"

$ws.Range("B4").Value = 19
$ws.Range("C4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("D4").Value = "Jewel"
$ws.Range("G4").Value = "Black or African American"

# --- Rows 10 / 11 / 12 rotate (Shadaisia, Kellie, Shaniek) ---
$ws.Range("B10").Value = 33
$ws.Range("C10").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("D10").Value = "Shaniek"

$ws.Range("B11").Value = 30
$ws.Range("C11").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("D11").Value = "Shadaisia"
$ws.Range("G11").Value = "Black or African American"

$ws.Range("B12").Value = 32
$ws.Range("C12").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("D12").Value = "Kellie"
$ws.Range("G12").Value = "White"

# --- Updated realeffort (F) scores for all 12 ranked rows ---
$ws.Range("F2").Value = 7.311265211180753
$ws.Range("F3").Value = 6.075952185643782
$ws.Range("F4").Value = 6.068676626552405
$ws.Range("F5").Value = 5.477047804629725
$ws.Range("F6").Value = 5.249471932023906
$ws.Range("F7").Value = 4.260356005502568
$ws.Range("F8").Value = 1.260598627945096
$ws.Range("F9").Value = 1.185192640848691
$ws.Range("F10").Value = 0.3346982378612178
$ws.Range("F11").Value = 0.26099946291021
$ws.Range("F12").Value = 0.07698541627100014
$ws.Range("F13").Value = 0.04919117767745862
